# Auto-generated: update Price (D) and Volume(1h) (E) columns
# per the commit diff ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.985.69"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "2.649.38"
$ws.Range("E3").Value = "  +5.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'113.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.46%  "
$ws.Range("D6").Value = "'326.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("D7").Value = "'0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.559"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("D10").Value = "'41.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").Value = "'20.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'0.0825"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "'7.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.12%  "
$ws.Range("D15").Value = "3.058.40"
$ws.Range("E15").Value = "  +5.71%  "
$ws.Range("D16").Value = "2.631.72"
$ws.Range("E16").Value = "  +5.09%  "
$ws.Range("D17").Value = "'0.874"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("D18").Value = "49.898.42"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").Value = "'13.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").Value = "'6.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").Value = "'72.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "'276.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "'26.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").Value = "'36.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").Value = "'5.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("E34").Value = "  +5.74%  "
$ws.Range("D35").Value = "'19.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'5.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.99%  "
$ws.Range("D38").Value = "'2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.22%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.71%  "
$ws.Range("D40").Value = "'124.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'21.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("D45").Value = "2.087.86"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("D46").Value = "'3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.95%  "
$ws.Range("D47").Value = "'2.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.46%  "
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").Value = "'9.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "'5.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").Value = "'59.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.81%  "
